$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4375.125
$ws.Range("I17").Value = 3000.25
$ws.Range("K17").Value = 9000.75
$ws.Range("M17").Value = -8832.75
$ws.Range("H18").Value = 4107.143
$ws.Range("I18").Value = 3250
$ws.Range("J18").Value = 5250
$ws.Range("K18").Value = 3250
$ws.Range("L18").Value = 5250
$ws.Range("M18").Value = -2966
$ws.Range("N18").Value = -5818
$ws.Range("H107").Value = 2007.1538
$ws.Range("I107").Value = 566.44446
$ws.Range("J107").Value = 5248.75
$ws.Range("K107").Value = 566.44446
$ws.Range("L107").Value = 5248.75
$ws.Range("M107").Value = 1353.55554
$ws.Range("N107").Value = -9088.75
$ws.Range("H113").Value = 4050.6
$ws.Range("J113").Value = 4891.75
$ws.Range("L113").Value = 4891.75
$ws.Range("N113").Value = -11399.75
$ws.Range("H132").Value = 4067.1667
$ws.Range("J132").Value = 999
$ws.Range("L132").Value = 2997
$ws.Range("N132").Value = -8057

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 877.8889
$ws.Range("I2").Value = 877.8889
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 877.8889
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -764.8889
$ws.Range("N2").ClearContents() | Out-Null
$ws.Range("H45").Value = 1699.75
$ws.Range("I45").Value = 1643.6666
$ws.Range("K45").Value = 1643.6666
$ws.Range("M45").Value = -1266.6666
$ws.Range("H88").Value = 3098
$ws.Range("I88").Value = 1424.3334
$ws.Range("K88").Value = 1424.3334
$ws.Range("M88").Value = -1018.3334
$ws.Range("H91").Value = 3098
$ws.Range("I91").Value = 1424.3334
$ws.Range("K91").Value = 1424.3334
$ws.Range("M91").Value = -20.33339999999998
$ws.Range("H116").Value = 877.8889
$ws.Range("I116").Value = 877.8889
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 877.8889
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1416.1111
$ws.Range("N116").ClearContents() | Out-Null
$ws.Range("H122").Value = 2559.375
$ws.Range("I122").Value = 2559.375
$ws.Range("K122").Value = 7678.125
$ws.Range("M122").Value = -5228.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 877.8889
$ws.Range("I3").Value = 877.8889
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 877.8889
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -763.8889
$ws.Range("N3").ClearContents() | Out-Null
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents() | Out-Null
$ws.Range("H86").Value = 1359.0834
$ws.Range("I86").Value = 1611.4
$ws.Range("J86").Value = 938.55554
$ws.Range("K86").Value = 1611.4
$ws.Range("L86").Value = 938.55554
$ws.Range("M86").Value = -488.4000000000001
$ws.Range("N86").Value = -3184.55554
$ws.Range("H89").Value = 1359.0834
$ws.Range("I89").Value = 1611.4
$ws.Range("J89").Value = 938.55554
$ws.Range("K89").Value = 8057
$ws.Range("L89").Value = 4692.7777
$ws.Range("M89").Value = -2441
$ws.Range("N89").Value = -15924.7777
$ws.Range("H99").Value = 1456.8572
$ws.Range("I99").Value = 1456.8572
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1456.8572
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 41.14280000000008
$ws.Range("N99").ClearContents() | Out-Null
$ws.Range("H105").Value = 423.5
$ws.Range("I105").Value = 418.8
$ws.Range("J105").Value = 447
$ws.Range("K105").Value = 418.8
$ws.Range("L105").Value = 447
$ws.Range("M105").Value = 1328.2
$ws.Range("N105").Value = -3941
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents() | Out-Null
$ws.Range("H134").Value = 8618.691999999999
$ws.Range("I134").Value = 8931
$ws.Range("J134").Value = 7916
$ws.Range("K134").Value = 26793
$ws.Range("L134").Value = 23748
$ws.Range("M134").Value = -24258
$ws.Range("N134").Value = -28818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 9700.5
$ws.Range("I32").Value = 10000
$ws.Range("J32").Value = 9401
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 9401
$ws.Range("M32").Value = -9684
$ws.Range("N32").Value = -10033
$ws.Range("H107").Value = 1096.375
$ws.Range("I107").Value = 1147
$ws.Range("K107").Value = 1147
$ws.Range("M107").Value = 773

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 251.25
$ws.Range("I12").Value = 186.66667
$ws.Range("J12").Value = 272.77777
$ws.Range("K12").Value = 560.00001
$ws.Range("L12").Value = 818.33331
$ws.Range("M12").Value = -387.00001
$ws.Range("N12").Value = -1164.33331
$ws.Range("H37").Value = 98750
$ws.Range("J37").Value = 98750
$ws.Range("L37").Value = 296250
$ws.Range("N37").Value = -296474
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents() | Out-Null
$ws.Range("N99").ClearContents() | Out-Null
$ws.Range("H103").Value = 569.2
$ws.Range("I103").Value = 569.2
$ws.Range("K103").Value = 1707.6
$ws.Range("M103").Value = -828.6000000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1750
$ws.Range("I113").Value = 1700
$ws.Range("J113").Value = 1800
$ws.Range("K113").Value = 1700
$ws.Range("L113").Value = 1800
$ws.Range("M113").Value = 470
$ws.Range("N113").Value = -6140
$ws.Range("H122").Value = 8958415
$ws.Range("I122").Value = 11400711
$ws.Range("J122").Value = 3331
$ws.Range("K122").Value = 34202133
$ws.Range("L122").Value = 9993
$ws.Range("M122").Value = -34199683
$ws.Range("N122").Value = -14893

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents() | Out-Null
$ws.Range("N25").ClearContents() | Out-Null
$ws.Range("H107").Value = 580
$ws.Range("I107").Value = 580
$ws.Range("K107").Value = 1740
$ws.Range("M107").Value = 180
